$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 54 (pushes the existing rows 54-59 down to 56-61),
# mirroring the weekly price-logging pattern described in the commit message
# ("Fruta / hortaliza, semanal").
$ws.Rows.Item(54).Insert()
$ws.Rows.Item(54).Insert()

# New row 54: Acelga, Primera, fecha 2022-03-17 (serial 44637)
$ws.Range("A54").Value = 1
$ws.Range("B54").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C54").Value = "Arica y Parinacota"
$ws.Range("D54").Value = 44637
$ws.Range("E54").Value = 15
$ws.Range("F54").Value = 100112009
$ws.Range("G54").Value = "Acelga"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 200
$ws.Range("K54").Value = 1400
$ws.Range("L54").Value = 1500
$ws.Range("M54").Value = 1450
$ws.Range("N54").Value = "$/atado 2,5 a 3 kilos"
$ws.Range("O54").Value = "Región de Arica y Parinacota"
$ws.Range("P54").Value = 483
$ws.Range("Q54").Value = 3
$ws.Range("R54").Value = "Hortaliza"

# New row 55: Acelga, Segunda, fecha 2022-03-17 (serial 44637)
$ws.Range("A55").Value = 1
$ws.Range("B55").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C55").Value = "Arica y Parinacota"
$ws.Range("D55").Value = 44637
$ws.Range("E55").Value = 15
$ws.Range("F55").Value = 100112009
$ws.Range("G55").Value = "Acelga"
$ws.Range("H55").Value = "Sin especificar"
$ws.Range("I55").Value = "Segunda"
$ws.Range("J55").Value = 160
$ws.Range("K55").Value = 1000
$ws.Range("L55").Value = 1200
$ws.Range("M55").Value = 1100
$ws.Range("N55").Value = "$/atado 2,5 a 3 kilos"
$ws.Range("O55").Value = "Región de Arica y Parinacota"
$ws.Range("P55").Value = 367
$ws.Range("Q55").Value = 3
$ws.Range("R55").Value = "Hortaliza"

# Ensure the date cells keep the same date/time number format used by the
# rest of the "Fecha" column (style index 2 in the original workbook).
$ws.Range("D54").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D55").NumberFormat = "YYYY-MM-DD HH:MM:SS"
